# Daily attendance processing - 2025-12-31 09:59:24
# Re-order the comma-separated "Recorded By" list in column G: rotate the
# list left by one position (move the first entry to the end of the list).
# This affects every data row on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colRecordedBy = 7  # Column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colRecordedBy)
    $val = $cell.Value2
    if ($null -eq $val) {
        continue
    }
    $strVal = [string]$val
    if ($strVal.IndexOf(",") -lt 0) {
        # Only one entry - rotation has no effect.
        continue
    }

    $rawParts = $strVal.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    # Rotate left by one: move the first element to the end of the list.
    $newParts = @()
    for ($i = 1; $i -lt $parts.Count; $i++) {
        $newParts += $parts[$i]
    }
    $newParts += $parts[0]

    $newVal = [string]::Join(", ", $newParts)
    $cell.Value2 = $newVal
}
